# Actualizacion desde MV -datos-
# Append the new daily "Diaria" rows (20-09-2021 .. 01-10-2021) below the
# existing data table, following the same layout as the prior rows:
#   A: Serie (date, stored as text)
#   B: Cupo (millones de pesos)
#   C: Monto demandado (millones de pesos)
#   D: Total monto adjudicado (millones de pesos)
#   E: Monto adjudicado bancos y sociedades financieras (millones de pesos)
#   F: Monto adjudicado AFP y otros (millones de pesos)
#   G: Tasas de interes base 360 dias (porcentaje)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Serie (date text), Cupo, Monto demandado,
# Total monto adjudicado, Monto adjudicado bancos, Monto adjudicado AFP, Tasa
$rows = @(
    @(58, "20-09-2021", 50000,  95000, 50000, 45000,  5000, 2),
    @(59, "21-09-2021", 50000, $null,     0, $null, $null, $null),
    @(60, "22-09-2021", 50000, $null,     0, $null, $null, $null),
    @(61, "23-09-2021", 50000, $null,     0, $null, $null, $null),
    @(62, "24-09-2021", 50000, 110000, 25000, 25000,     0, 2.1),
    @(63, "27-09-2021", 50000, $null,     0, $null, $null, $null),
    @(64, "28-09-2021", 50000, 100000, 50000, 40000, 10000, 2.18),
    @(65, "29-09-2021", 50000,  75000, 50000, 40000, 10000, 2.2),
    @(66, "30-09-2021", 50000, $null,     0, $null, $null, $null),
    @(67, "01-10-2021", 50000, $null,     0, $null, $null, $null)
)

foreach ($r in $rows) {
    $rowIndex = $r[0]
    $serie = $r[1]
    $cellA = $ws.Cells.Item($rowIndex, 1)

    # Some "dd-mm-yyyy" "Serie" strings (e.g. "01-10-2021") are also
    # unambiguous "mm-dd-yyyy" dates, so a plain .Value2 assignment gets
    # silently reinterpreted as a date serial by the COM layer. Writing the
    # literal as a quoted text formula and then converting it to a static
    # value (copy / paste-values) keeps every Serie entry a plain text
    # string - consistent with the rest of the column - without leaving any
    # number-format/style residue behind.
    $cellA.Formula = '="' + $serie + '"'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163) # xlPasteValues
    $excel.CutCopyMode = $false

    $ws.Cells.Item($rowIndex, 2).Value2 = $r[2]

    if ($null -ne $r[3]) {
        $ws.Cells.Item($rowIndex, 3).Value2 = $r[3]
    }

    $ws.Cells.Item($rowIndex, 4).Value2 = $r[4]

    if ($null -ne $r[5]) {
        $ws.Cells.Item($rowIndex, 5).Value2 = $r[5]
    }
    if ($null -ne $r[6]) {
        $ws.Cells.Item($rowIndex, 6).Value2 = $r[6]
    }
    if ($null -ne $r[7]) {
        $ws.Cells.Item($rowIndex, 7).Value2 = $r[7]
    }
}
